$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8684868812561035
$ws.Range("B1").Value = 1.323819637298584
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.387317895889282
